$wb = $excel.ActiveWorkbook

$changesBySheet = @{
    'Citywide Totals' = @{ "F2" = 91; "D3" = 132; "J3" = 226; "B6" = 368; "C6" = 471; "D6" = 407; "E6" = 461; "F6" = 515; "G6" = 431; "I6" = 495; "J6" = 409; "B7" = 493; "C7" = 625; "D7" = 636; "E7" = 684; "F7" = 747; "G7" = 660; "I7" = 826; "J7" = 775 }
    'Garfield Park' = @{ "C6" = 34; "I6" = 32; "C7" = 39; "I7" = 48 }
    'Grand Crossing' = @{ "G6" = 15; "G7" = 27 }
    'South Chicago' = @{ "B4" = 12; "B5" = 15 }
    'Englewood' = @{ "D3" = 17; "D7" = 45 }
    'South Shore' = @{ "D4" = 6; "D5" = 10 }
    'By Neighborhood' = @{ "D10" = 3; "F19" = 23; "D28" = 45; "C32" = 39; "I32" = 48; "G36" = 27; "D47" = 14; "J47" = 16; "E53" = 81; "F53" = 79; "J53" = 119; "C63" = 7; "B80" = 15; "D82" = 10; "C85" = 15; "B98" = 493; "C98" = 625; "D98" = 636; "E98" = 684; "F98" = 747; "G98" = 660; "I98" = 826; "J98" = 775 }
    'Loop' = @{ "F2" = 7; "E6" = 63; "F6" = 58; "J6" = 60; "E7" = 81; "F7" = 79; "J7" = 119 }
    'Avondale' = @{ "D5" = 3; "D6" = 3 }
    'New City' = @{ "C4" = 6; "C5" = 7 }
    'Chatham' = @{ "F5" = 16; "F6" = 23 }
    'United Center' = @{ "C4" = 12; "C5" = 15 }
    'Lake View' = @{ "J3" = 11; "D5" = 12; "D6" = 14; "J6" = 16 }
}

foreach ($sheetName in $changesBySheet.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $changesBySheet[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}

"Done applying " + $changesBySheet.Count + " sheet updates"